$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Marco"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Salvo"

$ws1.Range("A1").Value = "Instance"
$ws1.Range("B1").Value = "Time"
$ws1.Range("C1").Value = "Method"
$ws1.Range("A2").Value = "d493"
$ws1.Range("B2").Value = "64251ms"
$ws1.Range("C2").Value = "Original Method"

$ws2.Range("A1").Value = "Instance"
$ws2.Range("E1").Value = "Value"
$ws2.Range("A2").Value = "ali535"
$ws2.Range("C1").Value = "Initial value"
$ws2.Range("G2").Value = "With delta evaluetion"
$ws2.Range("D1").Value = "Time ms"
$ws2.Range("B1").Value = "Initial time ms"
$ws2.Range("G1").Value = "Method"
$ws2.Range("A3").Value = "d493"
$ws2.Range("G3").Value = "With delta evaluetion"

$ws2.Range("B2").Value = 9
$ws2.Range("C2").Value = 253127
$ws2.Range("D2").Value = 58548
$ws2.Range("E2").Value = 230423
$ws2.Range("F2").Formula = "=(E2-C2)/(D2-B2)"

$ws2.Range("B3").Value = 7
$ws2.Range("C3").Value = 41665
$ws2.Range("D3").Value = 39836
$ws2.Range("E3").Value = 38548
$ws2.Range("F3").Formula = "=(E3-C3)/(D3-B3)"

$ws2.Range("B2:E3").HorizontalAlignment = -4108
$ws2.Range("F2:F3").NumberFormat = "0%"
